# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets
# to reflect the latest scraped totals (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet updates ---
$wsExhibit.Range("F2").Value  = 6978
$wsExhibit.Range("F4").Value  = 60
$wsExhibit.Range("F5").Value  = 455
$wsExhibit.Range("F7").Value  = 6859
$wsExhibit.Range("F8").Value  = 74
$wsExhibit.Range("F9").Value  = 202
$wsExhibit.Range("F10").Value = 1287
$wsExhibit.Range("F11").Value = 21
$wsExhibit.Range("F13").Value = 409
$wsExhibit.Range("F15").Value = 17
$wsExhibit.Range("F16").Value = 416
$wsExhibit.Range("F17").Value = 48
$wsExhibit.Range("F18").Value = 41
$wsExhibit.Range("F19").Value = 17
$wsExhibit.Range("F20").Value = 5238
$wsExhibit.Range("F21").Value = 120
$wsExhibit.Range("F22").Value = 171
$wsExhibit.Range("F23").Value = 653
$wsExhibit.Range("F25").Value = 238

# --- 全部类型 sheet updates ---
$wsAll.Range("F2").Value  = 6979
$wsAll.Range("F4").Value  = 60
$wsAll.Range("F5").Value  = 455
$wsAll.Range("F6").Value  = 157
$wsAll.Range("F7").Value  = 6859
$wsAll.Range("F9").Value  = 202
$wsAll.Range("F10").Value = 1287
$wsAll.Range("F11").Value = 21
$wsAll.Range("F12").Value = 108
$wsAll.Range("F13").Value = 409
$wsAll.Range("F14").Value = 149
$wsAll.Range("F17").Value = 48
$wsAll.Range("F18").Value = 41
$wsAll.Range("F19").Value = 17
$wsAll.Range("F21").Value = 5238
$wsAll.Range("F23").Value = 120
$wsAll.Range("F24").Value = 171
$wsAll.Range("F25").Value = 653
$wsAll.Range("F26").Value = 216
$wsAll.Range("F27").Value = 0
